$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows appended to the NAV history table (2024-09-02 .. 2024-09-27).
# Format column A as text first so date-like strings are not auto-converted
# into Excel date serial numbers (matches the existing plain-text date cells).
$ws.Range("A634:A653").NumberFormat = "@"

# Row 634
$ws.Range("A634").Value = "2024-09-02"
$ws.Range("C634").Value = 1092.650024414062
$ws.Range("D634").Value = 715.0499877929688
$ws.Range("E634").Value = 251.3500061035156
$ws.Range("F634").Value = 490.5
$ws.Range("G634").Value = 1505.25
$ws.Range("H634").Value = 29322.30020141602
$ws.Range("I634").Value = 0
$ws.Range("J634").Value = 257.7896181747883

# Row 635
$ws.Range("A635").Value = "2024-09-03"
$ws.Range("C635").Value = 1085.099975585938
$ws.Range("D635").Value = 710.7999877929688
$ws.Range("E635").Value = 251
$ws.Range("F635").Value = 488.8500061035156
$ws.Range("G635").Value = 1509
$ws.Range("H635").Value = 29201.39971923828
$ws.Range("I635").Value = -0.004123158188384413
$ws.Range("J635").Value = 256.7267107997304

# Row 636
$ws.Range("A636").Value = "2024-09-04"
$ws.Range("C636").Value = 1080.449951171875
$ws.Range("D636").Value = 722.4000244140625
$ws.Range("E636").Value = 250.5
$ws.Range("F636").Value = 484.1499938964844
$ws.Range("G636").Value = 1488.099975585938
$ws.Range("H636").Value = 29200.89978027344
$ws.Range("I636").Value = -0.00001712037675078922
$ws.Range("J636").Value = 256.7223155417195

# Row 637
$ws.Range("A637").Value = "2024-09-05"
$ws.Range("C637").Value = 1069.150024414062
$ws.Range("D637").Value = 733.8499755859375
$ws.Range("E637").Value = 251.1499938964844
$ws.Range("F637").Value = 495.6499938964844
$ws.Range("G637").Value = 1447.599975585938
$ws.Range("H637").Value = 29262.3996887207
$ws.Range("I637").Value = 0.002106096350113556
$ws.Range("J637").Value = 257.2629974734746

# Row 638
$ws.Range("A638").Value = "2024-09-06"
$ws.Range("C638").Value = 1049.349975585938
$ws.Range("D638").Value = 718.9000244140625
$ws.Range("E638").Value = 247.8000030517578
$ws.Range("F638").Value = 483
$ws.Range("G638").Value = 1418.050048828125
$ws.Range("H638").Value = 28702.20024108887
$ws.Range("I638").Value = -0.01914400232349252
$ws.Range("J638").Value = 252.3379540520938

# Row 639
$ws.Range("A639").Value = "2024-09-09"
$ws.Range("C639").Value = 1038.699951171875
$ws.Range("D639").Value = 700.1500244140625
$ws.Range("E639").Value = 243.8999938964844
$ws.Range("F639").Value = 474.75
$ws.Range("G639").Value = 1411.849975585938
$ws.Range("H639").Value = 28242.6496887207
$ws.Range("I639").Value = -0.01601098691069303
$ws.Range("J639").Value = 248.2977743726946

# Row 640
$ws.Range("A640").Value = "2024-09-10"
$ws.Range("C640").Value = 1035.800048828125
$ws.Range("D640").Value = 713.4000244140625
$ws.Range("E640").Value = 248.25
$ws.Range("F640").Value = 478.7999877929688
$ws.Range("G640").Value = 1424.449951171875
$ws.Range("H640").Value = 28522.85046386719
$ws.Range("I640").Value = 0.009921192885042528
$ws.Range("J640").Value = 250.7611844851729

# Row 641
$ws.Range("A641").Value = "2024-09-11"
$ws.Range("C641").Value = 976.2999877929688
$ws.Range("D641").Value = 725.4000244140625
$ws.Range("E641").Value = 241.5500030517578
$ws.Range("F641").Value = 472.2000122070312
$ws.Range("G641").Value = 1399.599975585938
$ws.Range("H641").Value = 27922.25028991699
$ws.Range("I641").Value = -0.02105680758348599
$ws.Range("J641").Value = 245.4809544740616

# Row 642
$ws.Range("A642").Value = "2024-09-12"
$ws.Range("C642").Value = 986.1500244140625
$ws.Range("D642").Value = 726.0499877929688
$ws.Range("E642").Value = 246.1499938964844
$ws.Range("F642").Value = 479.8500061035156
$ws.Range("G642").Value = 1403.150024414062
$ws.Range("H642").Value = 28182.30001831055
$ws.Range("I642").Value = 0.009313351384414074
$ws.Range("J642").Value = 247.7672048612599

# Row 643
$ws.Range("A643").Value = "2024-09-13"
$ws.Range("C643").Value = 992.0999755859375
$ws.Range("D643").Value = 724.25
$ws.Range("E643").Value = 245.6499938964844
$ws.Range("F643").Value = 485.3999938964844
$ws.Range("G643").Value = 1410.949951171875
$ws.Range("H643").Value = 28258.59951782227
$ws.Range("I643").Value = 0.002707355306775728
$ws.Range("J643").Value = 248.437998718186

# Row 644
$ws.Range("A644").Value = "2024-09-16"
$ws.Range("C644").Value = 988.4000244140625
$ws.Range("D644").Value = 733.6500244140625
$ws.Range("E644").Value = 243.8000030517578
$ws.Range("F644").Value = 489.9500122070312
$ws.Range("G644").Value = 1404.550048828125
$ws.Range("H644").Value = 28313.45072937012
$ws.Range("I644").Value = 0.00194104493795801
$ws.Range("J644").Value = 248.9202280379943

# Row 645
$ws.Range("A645").Value = "2024-09-17"
$ws.Range("C645").Value = 974.9500122070312
$ws.Range("D645").Value = 745.4000244140625
$ws.Range("E645").Value = 240.8000030517578
$ws.Range("F645").Value = 482.2999877929688
$ws.Range("G645").Value = 1400.25
$ws.Range("H645").Value = 28196.30033874512
$ws.Range("I645").Value = -0.004137623200533361
$ws.Range("J645").Value = 247.8902899273823

# Row 646
$ws.Range("A646").Value = "2024-09-18"
$ws.Range("C646").Value = 962.0499877929688
$ws.Range("D646").Value = 717.5499877929688
$ws.Range("E646").Value = 235.9499969482422
$ws.Range("F646").Value = 471.75
$ws.Range("G646").Value = 1391.300048828125
$ws.Range("H646").Value = 27572.89979553223
$ws.Range("I646").Value = -0.02210930284198537
$ws.Range("J646").Value = 242.4096084357902

# Row 647
$ws.Range("A647").Value = "2024-09-19"
$ws.Range("C647").Value = 967
$ws.Range("D647").Value = 728.5
$ws.Range("E647").Value = 237.5500030517578
$ws.Range("F647").Value = 459.9500122070312
$ws.Range("G647").Value = 1374.150024414062
$ws.Range("H647").Value = 27641.0502166748
$ws.Range("I647").Value = 0.002471645044516532
$ws.Range("J647").Value = 243.0087589432238

# Row 648
$ws.Range("A648").Value = "2024-09-20"
$ws.Range("C648").Value = 970.8499755859375
$ws.Range("D648").Value = 748.3499755859375
$ws.Range("E648").Value = 237.8500061035156
$ws.Range("F648").Value = 466.2999877929688
$ws.Range("G648").Value = 1380.550048828125
$ws.Range("H648").Value = 27960.69967651367
$ws.Range("I648").Value = 0.01156430227264067
$ws.Range("J648").Value = 245.8189856865425

# Row 649
$ws.Range("A649").Value = "2024-09-23"
$ws.Range("C649").Value = 971.7999877929688
$ws.Range("D649").Value = 750.2000122070312
$ws.Range("E649").Value = 236.4499969482422
$ws.Range("F649").Value = 471.1499938964844
$ws.Range("G649").Value = 1375.400024414062
$ws.Range("H649").Value = 27984.94996643066
$ws.Range("I649").Value = 0.0008672991090191444
$ws.Range("J649").Value = 246.0321842738084

# Row 650
$ws.Range("A650").Value = "2024-09-24"
$ws.Range("C650").Value = 977.2999877929688
$ws.Range("D650").Value = 735.9000244140625
$ws.Range("E650").Value = 237.3000030517578
$ws.Range("F650").Value = 476.7000122070312
$ws.Range("G650").Value = 1363.699951171875
$ws.Range("H650").Value = 27912.20024108887
$ws.Range("I650").Value = -0.002599601765558408
$ws.Range("J650").Value = 245.392598573186

# Row 651
$ws.Range("A651").Value = "2024-09-25"
$ws.Range("C651").Value = 963.5999755859375
$ws.Range("D651").Value = 730.0499877929688
$ws.Range("E651").Value = 238.3500061035156
$ws.Range("F651").Value = 473.7000122070312
$ws.Range("G651").Value = 1365.400024414062
$ws.Range("H651").Value = 27741.79995727539
$ws.Range("I651").Value = -0.006104867489544392
$ws.Range("J651").Value = 243.8945092759817

# Row 652
$ws.Range("A652").Value = "2024-09-26"
$ws.Range("C652").Value = 993.1500244140625
$ws.Range("D652").Value = 744.0999755859375
$ws.Range("E652").Value = 241.1999969482422
$ws.Range("F652").Value = 471.75
$ws.Range("G652").Value = 1329.949951171875
$ws.Range("H652").Value = 28111.79975891113
$ws.Range("I652").Value = 0.01333726730801793
$ws.Range("J652").Value = 247.1473955411533

# Row 653
$ws.Range("A653").Value = "2024-09-27"
$ws.Range("C653").Value = 993
$ws.Range("D653").Value = 735.4500122070312
$ws.Range("E653").Value = 239.5500030517578
$ws.Range("F653").Value = 497.2999877929688
$ws.Range("G653").Value = 1392.199951171875
$ws.Range("H653").Value = 28306.40000915527
$ws.Range("I653").Value = 0.006922368966521059
$ws.Range("J653").Value = 248.8582410022039

